# Update cryptocurrency price/volume data and a few coin identity cells
# to reflect the latest scrape (GitHub Actions refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D2").Value = "64.008.74"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "3.160.30"
$ws.Range("E3").Value = "  -3.34%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "568.47"
$ws.Range("E5").Value = "  -2.30%  "
$ws.Range("D6").Value = "167.89"
$ws.Range("E6").Value = "  -5.90%  "
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  -5.20%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.158.96"
$ws.Range("E9").Value = "  -3.27%  "
$ws.Range("D10").Value = "0.120"
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("D11").Value = "6.69"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "0.383"
$ws.Range("E12").Value = "  -4.11%  "
$ws.Range("D13").Value = "3.710.45"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "64.098.75"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").Value = "25.23"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").Value = "0.0000158"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "3.157.23"
$ws.Range("E18").Value = "  -3.51%  "
$ws.Range("D19").Value = "414.48"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").Value = "12.78"
$ws.Range("E20").Value = "  -2.49%  "
$ws.Range("D21").Value = "5.33"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").Value = "7.08"
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "69.99"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "0.201"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").Value = "0.490"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").Value = "0.0000106"
$ws.Range("E27").Value = "  -4.94%  "
$ws.Range("D28").Value = "8.73"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.82"
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "21.71"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "5.01"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("D34").Value = "6.32"
$ws.Range("E34").Value = "  -3.43%  "
$ws.Range("D35").Value = "1.13"
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("D36").Value = "154.92"
$ws.Range("E36").Value = "  -2.52%  "
$ws.Range("D37").Value = "1.37"
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("D38").Value = "2.696.77"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").Value = "1.69"
$ws.Range("E39").Value = "  -5.05%  "
$ws.Range("D40").Value = "24.47"
$ws.Range("E40").Value = "  -6.49%  "
$ws.Range("D41").Value = "4.17"
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("D42").Value = "38.75"
$ws.Range("E42").Value = "  -2.71%  "
$ws.Range("D43").Value = "0.706"
$ws.Range("E43").Value = "  -7.47%  "
$ws.Range("D44").Value = "0.0623"
$ws.Range("E44").Value = "  -4.55%  "
$ws.Range("D45").Value = "5.62"
$ws.Range("E45").Value = "  -3.66%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "21.76"
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("D47").Value = "0.0261"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "294.62"
$ws.Range("E48").Value = "  -5.91%  "
$ws.Range("D49").Value = "2.04"
$ws.Range("E49").Value = "  -9.84%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "0.0990"
$ws.Range("E51").Value = "  -4.12%  "
